$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "Muskan Vaswan"
$ws.Range("B39").Value = "2021-01-15 11:02:48.019338"

$ws.Range("A40").Value = "Muskan Vaswan"
$ws.Range("B40").Value = "2021-01-16 18:40:41.957364"
